$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Update results for Steel": the Hydrogen-demand result for the
# "Iron & steel" industry (cell B3) is refreshed with the new figure.
$ws.Range("B3").Value = 497.5207405739714
